$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '24.892.69'
Set-TextCell 2 5 '  +0.46%  '

# Row 3
Set-TextCell 3 4 '1.707.42'
Set-TextCell 3 5 '  +0.37%  '

# Row 4
Set-TextCell 4 4 '0.9979'
Set-TextCell 4 5 '  -0.62%  '

# Row 5
Set-TextCell 5 4 '317.38'
Set-TextCell 5 5 '  +0.09%  '

# Row 6
Set-TextCell 6 4 '0.9990'
Set-TextCell 6 5 '  -0.48%  '

# Row 7
Set-TextCell 7 4 '0.3945'
Set-TextCell 7 5 '  -0.27%  '

# Row 8
Set-TextCell 8 5 '  -0.63%  '

# Row 9
Set-TextCell 9 4 '1.489'
Set-TextCell 9 5 '  -1.11%  '

# Row 10
Set-TextCell 10 4 '0.9986'
Set-TextCell 10 5 '  -0.59%  '

# Row 11
Set-TextCell 11 4 '53.24'
Set-TextCell 11 5 '  +0.69%  '

# Row 12
Set-TextCell 12 4 '0.08813'
Set-TextCell 12 5 '  -1.10%  '

# Row 13
Set-TextCell 13 4 '26.54'
Set-TextCell 13 5 '  +9.42%  '

# Row 14
Set-TextCell 14 4 '7.500'
Set-TextCell 14 5 '  -2.75%  '

# Row 15
Set-TextCell 15 4 '8.144'
Set-TextCell 15 5 '  -0.31%  '

# Row 16
Set-TextCell 16 4 '0.00001359'
Set-TextCell 16 5 '  +2.05%  '

# Row 17
Set-TextCell 17 4 '1.743.37'
Set-TextCell 17 5 '  +1.97%  '

# Row 18
Set-TextCell 18 4 '96.49'
Set-TextCell 18 5 '  -3.33%  '

# Row 19
Set-TextCell 19 4 '0.07180'
Set-TextCell 19 5 '  +0.35%  '

# Row 20
Set-TextCell 20 4 '21.10'
Set-TextCell 20 5 '  +5.09%  '

# Row 21
Set-TextCell 21 4 '7.292'
Set-TextCell 21 5 '  +0.84%  '

# Row 22
Set-TextCell 22 4 '0.9994'
Set-TextCell 22 5 '  -0.77%  '

# Row 23
Set-TextCell 23 4 '14.35'
Set-TextCell 23 5 '  -2.35%  '

# Row 24
Set-TextCell 24 4 '24.881.86'
Set-TextCell 24 5 '  +0.44%  '

# Row 25
Set-TextCell 25 4 '2.995'
Set-TextCell 25 5 '  -3.45%  '

# Row 26
Set-TextCell 26 4 '2.336'
Set-TextCell 26 5 '  -0.09%  '

# Row 27
Set-TextCell 27 4 '23.21'
Set-TextCell 27 5 '  +0.55%  '

# Row 28
Set-TextCell 28 2 'HuobiToken'
Set-TextCell 28 3 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 28 4 '6.148'
Set-TextCell 28 5 '  +18.29%  '

# Row 29
Set-TextCell 29 2 'Monero'
Set-TextCell 29 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 29 4 '166.43'
Set-TextCell 29 5 '  +0.93%  '

# Row 30
Set-TextCell 30 4 '145.42'
Set-TextCell 30 5 '  +4.38%  '

# Row 31
Set-TextCell 31 4 '8.456'
Set-TextCell 31 5 '  -9.28%  '

# Row 32
Set-TextCell 32 2 'WrappedliquidstakedEther2.0'
Set-TextCell 32 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 32 4 '1.935.06'
Set-TextCell 32 5 '  +1.92%  '

# Row 33
Set-TextCell 33 2 'WEMIXTOKEN'
Set-TextCell 33 3 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 33 4 '2.248'
Set-TextCell 33 5 '  +14.30%  '

# Row 34
Set-TextCell 34 4 '0.08818'
Set-TextCell 34 5 '  -3.90%  '

# Row 35
Set-TextCell 35 4 '0.03208'
Set-TextCell 35 5 '  +5.04%  '

# Row 36
Set-TextCell 36 2 'InternetComputer(DFINITY)'
Set-TextCell 36 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 36 4 '7.188'
Set-TextCell 36 5 '  -11.56%  '

# Row 37
Set-TextCell 37 2 'ImmutableX'
Set-TextCell 37 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 37 4 '1.041'
Set-TextCell 37 5 '  -3.63%  '

# Row 38
Set-TextCell 38 4 '0.2877'
Set-TextCell 38 5 '  +2.07%  '

# Row 39
Set-TextCell 39 4 '10.95'
Set-TextCell 39 5 '  -1.37%  '

# Row 40
Set-TextCell 40 4 '0.8358'
Set-TextCell 40 5 '  +6.81%  '

# Row 41
Set-TextCell 41 4 '0.09243'

# Row 42
Set-TextCell 42 4 '14.11'
Set-TextCell 42 5 '  -3.10%  '

# Row 43
Set-TextCell 43 4 '1.476'
Set-TextCell 43 5 '  +0.16%  '

# Row 44
Set-TextCell 44 4 '17.41'
Set-TextCell 44 5 '  +6.66%  '

# Row 45
Set-TextCell 45 4 '2.686'
Set-TextCell 45 5 '  +1.76%  '

# Row 46
Set-TextCell 46 4 '0.7401'
Set-TextCell 46 5 '  +1.91%  '

# Row 47
Set-TextCell 47 4 '4.244'
Set-TextCell 47 5 '  -0.15%  '

# Row 48
Set-TextCell 48 4 '1.394'
Set-TextCell 48 5 '  +2.37%  '

# Row 49
Set-TextCell 49 4 '0.9988'
Set-TextCell 49 5 '  -0.52%  '

# Row 50
Set-TextCell 50 4 '141.02'
Set-TextCell 50 5 '  +0.00%  '

# Row 51
Set-TextCell 51 4 '0.08322'
Set-TextCell 51 5 '  +3.23%  '
